# Weekly forward look stats update 25.07.25
#
# Changes:
#  1. Update the "as at <date>" sentence in A2 (17 July 2025 -> 25 July 2025).
#  2. Insert a new row at row 14 for the newly pre-announced publication
#     "Standard Determinate Sentence (SDS40) release data" (week 31, same
#     publication date as the rest of the 28 Jul 2025 week), which pushes
#     every following row down by one.
#  3. Correct the "Week" value for the "Knife and Offensive Weapon
#     Sentencing Statistics: January to March 2025" row (now row 17) from
#     34 to 33.
#  4. Extend the two conditional-formatting blocks so they keep covering
#     the whole table (now A5:F44 / A5:A44 instead of A5:F43 / A5:A43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "as at" sentence.
$ws.Range("A2").Value = "This list contains a week-by-week view of  MoJ Official and National Statistics that have been pre-announced on the gov.uk release calendar as at 25 July 2025"

# 2. Insert the new row for the SDS40 release data publication.
$ws.Rows(14).Insert()

$ws.Range("A14").Value = "28 Jul 2025"
$ws.Range("B14").Value = "Standard Determinate Sentence (SDS40) release data"
$ws.Range("C14").Value = "31 July 2025"
$ws.Range("D14").Value = "provisional"
$ws.Range("E14").Value = 31
$ws.Range("F14").Value = "standard"

# 3. Fix the "Week" number for the Knife and Offensive Weapon Sentencing
#    Statistics (January to March 2025) row, which is now row 17.
$ws.Range("E17").Value = 33

# 4. Extend the conditional formatting ranges to include the new last row.
$fc = $ws.Range("A5:F43").FormatConditions
$fc.Item(1).ModifyAppliesToRange($ws.Range("A5:F44"))
$fc.Item(4).ModifyAppliesToRange($ws.Range("A5:A44"))
